$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all merges first so cell writes below aren't blocked by
# Excel's 'can't change part of a merged cell' restriction.
$ws.Cells.UnMerge()

# Insert two new rows at row 154 (pushes old row 154+ down by two),
# matching the dimension growing from I211 to I213.
$ws.Rows.Item(154).Insert()
$ws.Rows.Item(154).Insert()

$data = @{}
$data[149] = @('Processes machinery waste', 'How waste will be managed on the site', 'Site activity details', '', '', '', 'Description of activities, processes, and end products including site operations, plant, ventilation, and machinery', 'string', 'MUST')
$data[150] = @('', '', 'Proposal waste management', '', '', '', 'Whether the proposal involves waste management development', 'boolean', 'MUST')
$data[151] = @('', '', 'Waste management[]', 'Waste management facility type', '', '', 'Type of waste management facility', 'enum', 'MUST')
$data[152] = @('', '', 'Waste management[]', 'Not applicable', '', '', 'Whether the facility is not applicable', 'boolean', 'MAY')
$data[153] = @('', '', 'Waste management[]', 'Total capacity', '', '', 'Total capacity of void in cubic metres (or tonnes/litres)', 'number', 'MAY')
$data[154] = @('', '', 'Waste management[]', 'Unit type', '', '', 'Unit for capacity/throughput (e.g. cubic metres, tonnes, litres)', 'enum', 'MAY')
$data[155] = @('', '', 'Waste management[]', 'Annual throughput', '', '', 'Maximum annual operational throughput in tonnes/litres', 'number', 'MAY')
$data[156] = @('', '', 'Waste management[]', 'Unit type', '', '', 'Unit for capacity/throughput (e.g. cubic metres, tonnes, litres)', 'enum', 'MAY')
$data[157] = @('', '', 'Waste streams throughput', 'Municipal', '', '', 'Maximum throughput for municipal waste (annual throughput in tonnes/litres)', 'number', 'MAY')
$data[158] = @('', '', 'Waste streams throughput', 'Construction demolition', '', '', 'Maximum throughput for construction and demolition waste (annual throughput in tonnes/litres)', 'number', 'MAY')
$data[159] = @('', '', 'Waste streams throughput', 'Commercial industrial', '', '', 'Maximum throughput for commercial and industrial waste (annual throughput in tonnes/litres)', 'number', 'MAY')
$data[160] = @('', '', 'Waste streams throughput', 'Hazardous', '', '', 'Maximum throughput for hazardous waste (annual throughput in tonnes/litres)', 'number', 'MAY')
$data[161] = @('Description of the proposal', 'What development, works or change of use is proposed', 'Proposal description', '', '', '', 'A description of what is being proposed, including the development, works, or change of use', 'string', 'MUST')
$data[162] = @('', '', 'Proposal started', '', '', '', 'Has any work on the proposal already been started', 'boolean', 'MUST')
$data[163] = @('', '', 'Proposal start date', '', '', '', 'The date when work on the proposal started, in YYYY-MM-DD format', 'date', 'MAY')
$data[164] = @('', '', 'Proposal completed', '', '', '', 'Has any work on the proposal already been completed', 'boolean', 'MUST')
$data[165] = @('', '', 'Proposal completion date', '', '', '', 'The date when work on the proposal was completed, in YYYY-MM-DD format', 'date', 'MAY')
$data[166] = @('', '', 'PIP reference', '', '', '', 'Reference number for the Planning in Principle (PIP) application this relates to', 'string', 'MAY')
$data[167] = @('', '', 'Is public service infrastructure', '', '', '', 'For applications made on or after 1 August 2021, is the proposal for public service infrastructure development', 'boolean', 'MUST')
$data[168] = @('Residential units', 'Details of the residential units that make up both the current and proposed development.', 'Residential unit change', '', '', '', 'Proposal includes the gain, loss or change of use of residential units', 'boolean', 'MUST')
$data[169] = @('', '', 'Residential unit summary[]', 'Tenure type', '', '', 'Category of housing tenure', 'enum', 'MUST')
$data[170] = @('', '', 'Residential unit summary[]', 'Housing type', '', '', 'Type of housing', 'enum', 'MUST')
$data[171] = @('', '', 'Residential unit summary[]', 'Existing unit breakdown[]', 'Units unknown', '', 'Whether the number of units is unknown', 'boolean', 'MUST')
$data[172] = @('', '', 'Residential unit summary[]', 'Existing unit breakdown[]', 'Units per bedroom number[]', 'No bedrooms unknown', 'Set to true when counting units where bedroom number is unknown', 'boolean', 'MUST')
$data[173] = @('', '', 'Residential unit summary[]', 'Existing unit breakdown[]', 'Units per bedroom number[]', 'Number of bedrooms', 'The number of bedrooms in unit', 'number', 'MAY')
$data[174] = @('', '', 'Residential unit summary[]', 'Existing unit breakdown[]', 'Units per bedroom number[]', 'Units', 'The number of units of that bedroom count', 'number', 'MUST')
$data[175] = @('', '', 'Residential unit summary[]', 'Existing unit breakdown[]', 'Total units', '', 'Total number of units', 'number', 'MAY')
$data[176] = @('', '', 'Residential unit summary[]', 'Proposed unit breakdown[]', 'Units unknown', '', 'Whether the number of units is unknown', 'boolean', 'MUST')
$data[177] = @('', '', 'Residential unit summary[]', 'Proposed unit breakdown[]', 'Units per bedroom number[]', 'No bedrooms unknown', 'Set to true when counting units where bedroom number is unknown', 'boolean', 'MUST')
$data[178] = @('', '', 'Residential unit summary[]', 'Proposed unit breakdown[]', 'Units per bedroom number[]', 'Number of bedrooms', 'The number of bedrooms in unit', 'number', 'MAY')
$data[179] = @('', '', 'Residential unit summary[]', 'Proposed unit breakdown[]', 'Units per bedroom number[]', 'Units', 'The number of units of that bedroom count', 'number', 'MUST')
$data[180] = @('', '', 'Residential unit summary[]', 'Proposed unit breakdown[]', 'Total units', '', 'Total number of units', 'number', 'MAY')
$data[181] = @('', '', 'Total existing units', '', '', '', 'The total number of existing units', 'number', 'MUST')
$data[182] = @('', '', 'Total proposed units', '', '', '', 'The total number of proposed units', 'number', 'MUST')
$data[183] = @('', '', 'Net change', '', '', '', 'Calculated net change in units', 'number', 'MUST')
$data[184] = @('Site area', 'How big the site is including relevant measurements', 'Site area in hectares', '', '', '', 'The size of the site in hectares', 'number', 'MUST')
$data[185] = @('', '', 'Site area provided by', '', '', '', 'Who provided the site area value', 'enum', 'MAY')
$data[186] = @('Site details', 'Where the proposed development will be built.', 'Site locations[]', 'Site boundary', '', '', 'Geometry of the site of the development, typically in GeoJSON format', 'wkt', 'MAY')
$data[187] = @('', '', 'Site locations[]', 'Address Text', '', '', 'Flexible field for capturing addresses', 'string', 'MAY')
$data[188] = @('', '', 'Site locations[]', 'Postcode', '', '', 'The postal code', 'string', 'MAY')
$data[189] = @('', '', 'Site locations[]', 'Easting', '', '', 'Easting coordinate in British National Grid (EPSG:27700)', 'number', 'MAY')
$data[190] = @('', '', 'Site locations[]', 'Northing', '', '', 'Northing coordinate in British National Grid (EPSG:27700)', 'number', 'MAY')
$data[191] = @('', '', 'Site locations[]', 'Latitude', '', '', 'Latitude coordinate in WGS84 (EPSG:4326)', 'number', 'MAY')
$data[192] = @('', '', 'Site locations[]', 'Longitude', '', '', 'Longitude coordinate in WGS84 (EPSG:4326)', 'number', 'MAY')
$data[193] = @('', '', 'Site locations[]', 'Description', '', '', 'A text description providing details about the subject. For parking changes, this describes how the proposed works affect existing car parking arrangements.', 'string', 'MAY')
$data[194] = @('', '', 'Site locations[]', 'UPRNs[]', '', '', 'Unique Property Reference Numbers (UPRNs) for properties within the site boundary', 'string', 'MAY')
$data[195] = @('Site Visit Details', 'Information to help the planning authority arrange a site visit', 'Site seen from public area', '', '', '', 'Can site be seen from a public road, public footpath, bridleway or other public land', 'boolean', 'MUST')
$data[196] = @('', '', 'Site visit contact type', '', '', '', 'Indicates who the authority should contact to arrange a site visit', 'enum', 'MUST')
$data[197] = @('', '', 'Contact reference', '', '', '', 'The reference of the applicant or agent who should be contacted for site visits', 'string', 'MAY')
$data[198] = @('', '', 'Other site visit contact', 'Full name', '', '', 'The complete name of a person', 'string', 'MUST')
$data[199] = @('', '', 'Other site visit contact', 'Phone number', '', '', 'A phone number', 'string', 'MUST')
$data[200] = @('', '', 'Other site visit contact', 'Email', '', '', 'The email address that can be used for electronic correspondence with the individual', 'string', 'MUST')
$data[201] = @('Trade effluent', 'Details of any liquid waste produced by industial processes on the proposed site, and how it will be diposed of.', 'Disposal required', '', '', '', 'Does the proposal involve the disposal of trade effluents or waste (true/false)', 'boolean', 'MUST')
$data[202] = @('', '', 'Description', '', '', '', 'describe the nature, volume and means of disposal of trade effluents or waste', 'string', 'MAY')
$data[203] = @('Trees and hedges information', 'Details of trees and/or hedges that will be affected by the proposed development', 'Trees on site', '', '', '', 'Whether trees or hedges are present on the proposed development site', 'boolean', 'MUST')
$data[204] = @('', '', 'Trees on adjacent land', '', '', '', 'Whether trees or hedges on land adjacent to the proposed development site could influence the development or might be important as part of the local landscape character', 'boolean', 'MUST')
$data[205] = @('Vehicle parking', 'Details of current parking facilities at the site and any changes that would be made by the proposed development.', 'Parking spaces[]', 'Parking space type', '', '', 'Type of parking space or vehicle type', 'enum', 'MUST')
$data[206] = @('', '', 'Parking spaces[]', 'Vehicle type other', '', '', 'Vehicle type when parking space type is ''other''', 'string', 'MAY')
$data[207] = @('', '', 'Parking spaces[]', 'Total existing', '', '', 'Total number of existing parking spaces', 'number', 'MUST')
$data[208] = @('', '', 'Parking spaces[]', 'Total proposed', '', '', 'Total number of proposed parking spaces', 'number', 'MUST')
$data[209] = @('', '', 'Parking spaces[]', 'Difference in spaces', '', '', 'Net change in parking spaces (proposed minus existing)', 'number', 'MUST')
$data[210] = @('Waste storage and collection', 'Any waste storage or recycling arrangements are in place, such as waste storage areas', 'Needs waste storage area', '', '', '', 'Does the proposal require a waste storage area', 'boolean', 'MUST')
$data[211] = @('', '', 'Waste storage area details', '', '', '', 'Details of the waste storage area including location, size, design and access arrangements', 'string', 'MAY')
$data[212] = @('', '', 'Separate recycling arrangements', '', '', '', 'Does the proposal include separate recycling arrangements', 'boolean', 'MUST')
$data[213] = @('', '', 'Separate recycling arrangements details', '', '', '', 'Details of the recycling arrangements including types of materials, collection methods and storage facilities', 'string', 'MAY')

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}

# Re-apply the merged cell ranges for columns A and B.
$ws.Range("B65").Merge() | Out-Null
$ws.Range("B205:B209").Merge() | Out-Null
$ws.Range("B103:B109").Merge() | Out-Null
$ws.Range("A2:A19").Merge() | Out-Null
$ws.Range("B168:B183").Merge() | Out-Null
$ws.Range("B79:B88").Merge() | Out-Null
$ws.Range("A72:A78").Merge() | Out-Null
$ws.Range("A99:A102").Merge() | Out-Null
$ws.Range("A30:A37").Merge() | Out-Null
$ws.Range("A144:A148").Merge() | Out-Null
$ws.Range("B51:B64").Merge() | Out-Null
$ws.Range("A79:A88").Merge() | Out-Null
$ws.Range("B66:B68").Merge() | Out-Null
$ws.Range("A117:A128").Merge() | Out-Null
$ws.Range("B195:B200").Merge() | Out-Null
$ws.Range("B186:B194").Merge() | Out-Null
$ws.Range("B42:B47").Merge() | Out-Null
$ws.Range("A95:A98").Merge() | Out-Null
$ws.Range("B38:B41").Merge() | Out-Null
$ws.Range("A129:A143").Merge() | Out-Null
$ws.Range("B210:B213").Merge() | Out-Null
$ws.Range("B110:B116").Merge() | Out-Null
$ws.Range("A89:A94").Merge() | Out-Null
$ws.Range("B149:B160").Merge() | Out-Null
$ws.Range("A48:A50").Merge() | Out-Null
$ws.Range("B144:B148").Merge() | Out-Null
$ws.Range("A161:A167").Merge() | Out-Null
$ws.Range("A26:A29").Merge() | Out-Null
$ws.Range("A51:A64").Merge() | Out-Null
$ws.Range("A210:A213").Merge() | Out-Null
$ws.Range("B201:B202").Merge() | Out-Null
$ws.Range("A184:A185").Merge() | Out-Null
$ws.Range("A65").Merge() | Out-Null
$ws.Range("A203:A204").Merge() | Out-Null
$ws.Range("A69:A71").Merge() | Out-Null
$ws.Range("B99:B102").Merge() | Out-Null
$ws.Range("B117:B128").Merge() | Out-Null
$ws.Range("A186:A194").Merge() | Out-Null
$ws.Range("A195:A200").Merge() | Out-Null
$ws.Range("B72:B78").Merge() | Out-Null
$ws.Range("B95:B98").Merge() | Out-Null
$ws.Range("A66:A68").Merge() | Out-Null
$ws.Range("A110:A116").Merge() | Out-Null
$ws.Range("B89:B94").Merge() | Out-Null
$ws.Range("A149:A160").Merge() | Out-Null
$ws.Range("A42:A47").Merge() | Out-Null
$ws.Range("B48:B50").Merge() | Out-Null
$ws.Range("A38:A41").Merge() | Out-Null
$ws.Range("A205:A209").Merge() | Out-Null
$ws.Range("B161:B167").Merge() | Out-Null
$ws.Range("B26:B29").Merge() | Out-Null
$ws.Range("B30:B37").Merge() | Out-Null
$ws.Range("A20:A25").Merge() | Out-Null
$ws.Range("A103:A109").Merge() | Out-Null
$ws.Range("B129:B143").Merge() | Out-Null
$ws.Range("B20:B25").Merge() | Out-Null
$ws.Range("A168:A183").Merge() | Out-Null
$ws.Range("B2:B19").Merge() | Out-Null
$ws.Range("B69:B71").Merge() | Out-Null
$ws.Range("B184:B185").Merge() | Out-Null
$ws.Range("A201:A202").Merge() | Out-Null
$ws.Range("B203:B204").Merge() | Out-Null

$ws.Range("A1").Select()
